$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text looks like a number need the NumberFormat trick
# so Excel stores them as literal text (matching the inlineStr source),
# then we reset the style back to Normal to avoid leaving a stray format.
$textForceCells = @("D5", "D6", "D8", "D10", "D11", "D12", "D19", "D22", "D24", "D25", "D26", "D28", "D31", "D32", "D33", "D36", "D38", "D42", "D43", "D44", "D45", "D48", "D49", "D50")
foreach ($addr in $textForceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '67.789.90'
$ws.Range('E2').Value = '  +8.90%  '
$ws.Range('D3').Value = '3.515.90'
$ws.Range('E3').Value = '  +10.66%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = '191.31'
$ws.Range('E5').Value = '  +11.99%  '
$ws.Range('D6').Value = '555.68'
$ws.Range('E6').Value = '  +8.72%  '
$ws.Range('D7').Value = '3.506.60'
$ws.Range('E7').Value = '  +10.43%  '
$ws.Range('D8').Value = '0.612'
$ws.Range('E8').Value = '  +4.31%  '
$ws.Range('E9').Value = '  -0.07%  '
$ws.Range('D10').Value = '0.640'
$ws.Range('E10').Value = '  +8.55%  '
$ws.Range('D11').Value = '57.15'
$ws.Range('E11').Value = '  +5.53%  '
$ws.Range('D12').Value = '0.152'
$ws.Range('E12').Value = '  +17.40%  '
$ws.Range('E13').Value = '  +10.78%  '
$ws.Range('E14').Value = '  +7.92%  '
$ws.Range('D15').Value = '4.082.68'
$ws.Range('D16').Value = '3.522.86'
$ws.Range('E16').Value = '  +11.33%  '
$ws.Range('D17').Value = '68.142.41'
$ws.Range('E17').Value = '  +9.76%  '
$ws.Range('E18').Value = '  +7.27%  '
$ws.Range('D19').Value = '18.37'
$ws.Range('E19').Value = '  +8.71%  '
$ws.Range('E20').Value = '  +11.88%  '
$ws.Range('E21').Value = '  +7.96%  '
$ws.Range('D22').Value = '409.41'
$ws.Range('E22').Value = '  +13.63%  '
$ws.Range('E23').Value = '  +8.59%  '
$ws.Range('D24').Value = '84.85'
$ws.Range('E24').Value = '  +7.56%  '
$ws.Range('D25').Value = '11.68'
$ws.Range('E25').Value = '  +8.33%  '
$ws.Range('D26').Value = '4.21'
$ws.Range('E26').Value = '  +10.52%  '
$ws.Range('E27').Value = '  +12.37%  '
$ws.Range('D28').Value = '6.14'
$ws.Range('E28').Value = '  -0.26%  '
$ws.Range('E29').Value = '  +8.42%  '
$ws.Range('E30').Value = '  +7.52%  '
$ws.Range('D31').Value = '30.63'
$ws.Range('E31').Value = '  +10.24%  '
$ws.Range('D32').Value = '682.25'
$ws.Range('E32').Value = '  +10.19%  '
$ws.Range('D33').Value = '6.89'
$ws.Range('E33').Value = '  +7.88%  '
$ws.Range('E34').Value = '  +7.31%  '
$ws.Range('E35').Value = '  +9.47%  '
$ws.Range('D36').Value = '60.61'
$ws.Range('E36').Value = '  +7.06%  '
$ws.Range('D37').Value = '0.0₃0842'
$ws.Range('E37').Value = '  +27.25%  '
$ws.Range('D38').Value = '39.10'
$ws.Range('E38').Value = '  +8.41%  '
$ws.Range('E39').Value = '  +8.13%  '
$ws.Range('E40').Value = '  -0.20%  '
$ws.Range('E41').Value = '  +25.52%  '
$ws.Range('D42').Value = '0.134'
$ws.Range('E42').Value = '  +12.17%  '
$ws.Range('D43').Value = '2.75'
$ws.Range('E43').Value = '  +15.21%  '
$ws.Range('D44').Value = '3.03'
$ws.Range('E44').Value = '  +17.53%  '
$ws.Range('D45').Value = '1.00'
$ws.Range('E45').Value = '  +0.42%  '
$ws.Range('D46').Value = '3.046.64'
$ws.Range('E46').Value = '  +9.24%  '
$ws.Range('E47').Value = '  +11.26%  '
$ws.Range('B48').Value = 'WEMIXToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D48').Value = '2.75'
$ws.Range('E48').Value = '  +5.66%  '
$ws.Range('B49').Value = 'ApeXProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D49').Value = '3.23'
$ws.Range('E49').Value = '  +11.48%  '
$ws.Range('B50').Value = 'THORChain'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D50').Value = '9.08'
$ws.Range('E50').Value = '  +22.25%  '
$ws.Range('E51').Value = '  +8.04%  '

foreach ($addr in $textForceCells) {
    $ws.Range($addr).Style = "Normal"
}
